$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the existing four battle sheets.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Battle 1").Name = "Single STR"
$wb.Worksheets.Item("Battle 2").Name = "Single AGL"
$wb.Worksheets.Item("Battle 3").Name = "Shield"
$wb.Worksheets.Item("Battle 4").Name = "Group MANA - PC"

# ---------------------------------------------------------------------------
# 2) Add two new sheets at the end, using "Group MANA - PC" as a formatting
#    template (same column widths / layout), then wipe their contents so we
#    can rebuild them from scratch.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("Group MANA - PC")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$manaEnemy = $wb.Worksheets.Item($wb.Worksheets.Count)
$manaEnemy.Name = "Group MANA - Enemy"
$manaEnemy.Cells.Clear()

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet2)
$allEnemies = $wb.Worksheets.Item($wb.Worksheets.Count)
$allEnemies.Name = "All Enemies - PC"
$allEnemies.Cells.Clear()

# ---------------------------------------------------------------------------
# 3) Header row, identical on every battle sheet.
# ---------------------------------------------------------------------------
$headers = @("Index","NAME","ROLE","LIVES","Position","Initiative","CURRENT HP","CURRENT STR","CURRENT AGL","CURRENT MANA","CURRENT DEF","COMMAND","TARGET","Stoned","Cursed","Blinded","Stunned","Paralyzed","Poisoned","Confused","ACTIONS TAKEN")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $manaEnemy.Cells.Item(1, $i + 1).Value = $headers[$i]
    $allEnemies.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 4) "Group MANA - Enemy" sheet content.
#    Fill column-by-column (NAME first) so new shared-string entries land in
#    the same order the source workbook used: ROBO, Jerk, ColtGun.
# ---------------------------------------------------------------------------
$manaEnemy.Range("B2").Value = "HUME"
$manaEnemy.Range("B3").Value = "MUTE"
$manaEnemy.Range("B4").Value = "MONS"
$manaEnemy.Range("B5").Value = "ROBO"
$manaEnemy.Range("B6").Value = "Jerk"

$manaEnemy.Range("C2").Value = "Player"
$manaEnemy.Range("C3").Value = "Player"
$manaEnemy.Range("C4").Value = "Player"
$manaEnemy.Range("C5").Value = "Player"
$manaEnemy.Range("C6").Value = "Enemy"

$manaEnemy.Range("D2").Value = 1
$manaEnemy.Range("D3").Value = 1
$manaEnemy.Range("D4").Value = 1
$manaEnemy.Range("D5").Value = 1
$manaEnemy.Range("D6").Value = 1

$manaEnemy.Range("E2").Value = 1
$manaEnemy.Range("E3").Value = 2
$manaEnemy.Range("E4").Value = 3
$manaEnemy.Range("E5").Value = 4

$manaEnemy.Range("L2").Value = "LongSword"
$manaEnemy.Range("L3").Value = "Rapier"
$manaEnemy.Range("L4").Value = "Nail"
$manaEnemy.Range("L5").Value = "ColtGun"

$manaEnemy.Range("M2").Value = "Jerk"
$manaEnemy.Range("M3").Value = "Jerk"
$manaEnemy.Range("M4").Value = "Jerk"
$manaEnemy.Range("M5").Value = "Jerk"

# NAME formula column: A2/A3 individually, then A4:A6 as one filled block
# (matches the shared-formula grouping used elsewhere in the workbook).
$manaEnemy.Range("A2").Formula = "=B2"
$manaEnemy.Range("A3").Formula = "=B3"
$manaEnemy.Range("A4:A6").Formula = "=B4"

# ---------------------------------------------------------------------------
# 5) "All Enemies - PC" sheet content.
# ---------------------------------------------------------------------------
$allEnemies.Range("B2").Value = "Flammie"
$allEnemies.Range("C2").Value = "Player"
$allEnemies.Range("D2").Value = 1
$allEnemies.Range("E2").Value = 1
$allEnemies.Range("L2").Value = "Flame"

$allEnemies.Range("B3").Value = "Goblin"
$allEnemies.Range("C3").Value = "Enemy"
$allEnemies.Range("D3").Value = 4

$allEnemies.Range("B4").Value = "Jaguar"
$allEnemies.Range("C4").Value = "Enemy"
$allEnemies.Range("D4").Value = 4

$allEnemies.Range("B5").Value = "Eagle"
$allEnemies.Range("C5").Value = "Enemy"
$allEnemies.Range("D5").Value = 4

$allEnemies.Range("A2").Formula = "=B2"
$allEnemies.Range("A3").Formula = "=B3"
$allEnemies.Range("A4:A5").Formula = "=B4"

# ---------------------------------------------------------------------------
# 6) Selections / active tab: "All Enemies - PC" becomes the active sheet,
#    losing tabSelected on "Group MANA - PC" and gaining it on the new sheet.
# ---------------------------------------------------------------------------
$manaEnemy.Range("E9").Select()
$allEnemies.Range("A6").Select()
$allEnemies.Activate()
